$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 = "Save", matching the formatting of the existing
# header cells (e.g. G1 -> bold font, thin border, centered/top alignment)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new data cell H2 = 0 (numeric, unstyled like the other data cells)
$ws.Range("H2").Value = 0
